$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new "Jurisdiction" property row after "Contact" (row 10), before "Description" (old row 11)
$ws.Range("A11").EntireRow.Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update Version and Date metadata values
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
